$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2 = "31211092660604015457550760000146171184325883.xml"
$ws.Range("B2").Value2 = "C:\Users\Stars\Desktop\TRANSMOB\31211092660604015457550760000146171184325883.xml"
$ws.Range("A3").Value2 = "31211092660604015457550760000146221569723315.xml"
$ws.Range("B3").Value2 = "C:\Users\Stars\Desktop\TRANSMOB\31211092660604015457550760000146221569723315.xml"
$ws.Range("A4").Value2 = "31211092660604015457550760000146311031134436.xml"
$ws.Range("B4").Value2 = "C:\Users\Stars\Desktop\TRANSMOB\31211092660604015457550760000146311031134436.xml"
$ws.Range("A5").Value2 = "35210833931486000564550040000805401453515100.xml"
$ws.Range("B5").Value2 = "C:\Users\Stars\Desktop\TRANSMOB\35210833931486000564550040000805401453515100.xml"
$ws.Range("F5").Value2 = "True"
$ws.Range("A6").Value2 = "35210850505924000118570010002698471002698478-ret-cons-cte.xml"
$ws.Range("B6").Value2 = "C:\Users\Stars\Desktop\TRANSMOB\35210850505924000118570010002698471002698478-ret-cons-cte.xml"
$ws.Range("A7").Value2 = "35210850505924000118570010002729211002729218-ret-cons-cte.xml"
$ws.Range("B7").Value2 = "C:\Users\Stars\Desktop\TRANSMOB\35210850505924000118570010002729211002729218-ret-cons-cte.xml"
$ws.Range("E7").Value2 = "False"
$ws.Range("A8").Value2 = "35210950505924000118570010002738851002738853-ret-cons-cte.xml"
$ws.Range("B8").Value2 = "C:\Users\Stars\Desktop\TRANSMOB\35210950505924000118570010002738851002738853-ret-cons-cte.xml"
$ws.Range("E8").Value2 = "True"
$ws.Range("F8").Value2 = "True"
$ws.Range("A9").Value2 = "35211050505924000118570010002782611002782613-ret-cons-cte.xml"
$ws.Range("B9").Value2 = "C:\Users\Stars\Desktop\TRANSMOB\35211050505924000118570010002782611002782613-ret-cons-cte.xml"
$ws.Range("F9").Value2 = "True"
$ws.Range("A10").Value2 = "35211050505924000118570010002782621002782629-ret-cons-cte.xml"
$ws.Range("B10").Value2 = "C:\Users\Stars\Desktop\TRANSMOB\35211050505924000118570010002782621002782629-ret-cons-cte.xml"
$ws.Range("F10").Value2 = "True"
$ws.Range("A11").Value2 = "35211050505924000118570010002782671002782676-ret-cons-cte.xml"
$ws.Range("B11").Value2 = "C:\Users\Stars\Desktop\TRANSMOB\35211050505924000118570010002782671002782676-ret-cons-cte.xml"
$ws.Range("F11").Value2 = "True"
$ws.Range("A12").Value2 = "35211050505924000118570010002785081002785085-ret-cons-cte.xml"
$ws.Range("B12").Value2 = "C:\Users\Stars\Desktop\TRANSMOB\35211050505924000118570010002785081002785085-ret-cons-cte.xml"
$ws.Range("F12").Value2 = "False"
$ws.Range("A13").Value2 = "35211204400329000109570010000605251000605255-ret-cons-cte.xml"
$ws.Range("B13").Value2 = "C:\Users\Stars\Desktop\TRANSMOB\35211204400329000109570010000605251000605255-ret-cons-cte.xml"
$ws.Range("D13").Value2 = "False"
$ws.Range("A14").Value2 = "35211204400329000109570010000605261000605260-ret-cons-cte.xml"
$ws.Range("B14").Value2 = "C:\Users\Stars\Desktop\TRANSMOB\35211204400329000109570010000605261000605260-ret-cons-cte.xml"
$ws.Range("D14").Value2 = "False"
$ws.Range("A15").Value2 = "35211250505924000118570010002829241002829243-ret-cons-cte.xml"
$ws.Range("B15").Value2 = "C:\Users\Stars\Desktop\TRANSMOB\35211250505924000118570010002829241002829243-ret-cons-cte.xml"
$ws.Range("A16").Value2 = "35211250505924000118570010002829251002829259-ret-cons-cte.xml"
$ws.Range("B16").Value2 = "C:\Users\Stars\Desktop\TRANSMOB\35211250505924000118570010002829251002829259-ret-cons-cte.xml"
$ws.Range("A17").Value2 = "35211250505924000118570010002829261002829264-ret-cons-cte.xml"
$ws.Range("B17").Value2 = "C:\Users\Stars\Desktop\TRANSMOB\35211250505924000118570010002829261002829264-ret-cons-cte.xml"
$ws.Range("D17").Value2 = "True"
$ws.Range("A18").Value2 = "35211250505924000118570010002829291002829290-ret-cons-cte.xml"
$ws.Range("B18").Value2 = "C:\Users\Stars\Desktop\TRANSMOB\35211250505924000118570010002829291002829290-ret-cons-cte.xml"
$ws.Range("D18").Value2 = "True"
$ws.Range("A19").Value2 = "35211250505924000118570010002829411002829414-ret-cons-cte.xml"
$ws.Range("B19").Value2 = "C:\Users\Stars\Desktop\TRANSMOB\35211250505924000118570010002829411002829414-ret-cons-cte.xml"
$ws.Range("D19").Value2 = "True"
$ws.Range("A20").Value2 = "35211250505924000118570010002829461002829461-ret-cons-cte.xml"
$ws.Range("B20").Value2 = "C:\Users\Stars\Desktop\TRANSMOB\35211250505924000118570010002829461002829461-ret-cons-cte.xml"
$ws.Range("D20").Value2 = "True"
$ws.Range("A21").Value2 = "52211088305859002101550230000027711001975667.xml"
$ws.Range("B21").Value2 = "C:\Users\Stars\Desktop\TRANSMOB\52211088305859002101550230000027711001975667.xml"
$ws.Range("F21").Value2 = "False"
